$d = $word.ActiveDocument

$pairs = @(
    @("54÷8=", "76÷8="),
    @("20÷4=", "26÷9="),
    @("32÷6=", "64÷9="),
    @("63÷3=", "16÷2="),
    @("89÷7=", "68÷5="),
    @("87÷4=", "74÷5="),
    @("75÷3=", "88÷7="),
    @("33÷3=", "96÷4="),
    @("68÷7=", "69÷8="),
    @("13÷6=", "93÷5="),
    @("97÷3=", "61÷3="),
    @("70÷5=", "86÷4="),
    @("64÷4=", "85÷7="),
    @("77÷4=", "55÷5="),
    @("91÷5=", "45÷3="),
    @("56÷4=", "76÷2="),
    @("31÷2=", "21÷7="),
    @("41÷4=", "51÷5="),
    @("93÷2=", "63÷6="),
    @("87÷2=", "34÷6="),
    @("89÷8=", "34÷7="),
    @("60÷9=", "69÷4="),
    @("75÷9=", "81÷9="),
    @("42÷9=", "63÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
